$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 86, shifting existing rows 86..197 down to 87..198.
$ws.Rows.Item(86).Insert()

# Populate the newly inserted row 86 with the new data record.
$ws.Range("A86").Value = 5
$ws.Range("B86").Value = "Macroferia Regional de Talca"
$ws.Range("C86").Value = "Maule"
$ws.Range("D86").Value = 44483
$ws.Range("E86").Value = 7
$ws.Range("F86").Value = 100114013
$ws.Range("G86").Value = "Zanahoria"
$ws.Range("H86").Value = "Sin especificar"
$ws.Range("I86").Value = "Primera"
$ws.Range("J86").Value = 400
$ws.Range("K86").Value = 8000
$ws.Range("L86").Value = 8000
$ws.Range("M86").Value = 8000
$ws.Range("N86").Value = "$/saco 20 kilos"
$ws.Range("O86").Value = "Región de Ñuble"
$ws.Range("P86").Value = 400
$ws.Range("Q86").Value = 20
$ws.Range("R86").Value = "Hortaliza"
